# Scheduled runner update: refresh computed market-price / profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) on several leve rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 12072.956
$ws.Range("I62").Value = 18298.334
$ws.Range("J62").Value = 5281.636
$ws.Range("K62").Value = 18298.334
$ws.Range("L62").Value = 5281.636
$ws.Range("M62").Value = -17674.334
$ws.Range("N62").Value = -6529.636

$ws.Range("H65").Value = 12072.956
$ws.Range("I65").Value = 18298.334
$ws.Range("J65").Value = 5281.636
$ws.Range("K65").Value = 91491.67
$ws.Range("L65").Value = 26408.18
$ws.Range("M65").Value = -88371.67
$ws.Range("N65").Value = -32648.18

$ws.Range("H76").Value = 8397.808000000001
$ws.Range("I76").Value = 11024.929
$ws.Range("J76").Value = 5332.8335
$ws.Range("K76").Value = 11024.929
$ws.Range("L76").Value = 5332.8335
$ws.Range("M76").Value = -10709.929
$ws.Range("N76").Value = -5962.8335

$ws.Range("H79").Value = 8397.808000000001
$ws.Range("I79").Value = 11024.929
$ws.Range("J79").Value = 5332.8335
$ws.Range("K79").Value = 11024.929
$ws.Range("L79").Value = 5332.8335
$ws.Range("M79").Value = -9932.929
$ws.Range("N79").Value = -7516.8335

$ws.Range("H129").Value = 911.5
$ws.Range("I129").Value = 327
$ws.Range("J129").Value = 1496
$ws.Range("K129").Value = 981
$ws.Range("L129").Value = 4488
$ws.Range("M129").Value = 4019
$ws.Range("N129").Value = -14488

$ws.Range("H132").Value = 3031.0151
$ws.Range("I132").Value = 1372.6522
$ws.Range("K132").Value = 4117.9566
$ws.Range("M132").Value = -1587.9566

$ws.Range("H137").Value = 1743.0857
$ws.Range("I137").Value = 1548.85
$ws.Range("J137").Value = 2908.5
$ws.Range("K137").Value = 4646.549999999999
$ws.Range("L137").Value = 8725.5
$ws.Range("M137").Value = -2096.549999999999
$ws.Range("N137").Value = -13825.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1357.93
$ws.Range("I32").Value = 1210.3053
$ws.Range("J32").Value = 4162.8
$ws.Range("K32").Value = 1210.3053
$ws.Range("L32").Value = 4162.8
$ws.Range("M32").Value = -923.3053
$ws.Range("N32").Value = -4736.8

$ws.Range("H61").Value = 459436.25
$ws.Range("I61").Value = 368600.84
$ws.Range("K61").Value = 368600.84
$ws.Range("M61").Value = -368388.84

$ws.Range("H122").Value = 7331.4287
$ws.Range("I122").Value = 7280
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 21840
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -19390
$ws.Range("N122").Value = -28900

$ws.Range("H132").Value = 3154.6875
$ws.Range("I132").Value = 2896.6
$ws.Range("J132").Value = 3849.5386
$ws.Range("K132").Value = 8689.799999999999
$ws.Range("L132").Value = 11548.6158
$ws.Range("M132").Value = -6159.799999999999
$ws.Range("N132").Value = -16608.6158

$ws.Range("H136").Value = 459436.25
$ws.Range("I136").Value = 368600.84
$ws.Range("K136").Value = 1105802.52
$ws.Range("M136").Value = -1103252.52

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1499.4546
$ws.Range("I20").Value = 1378.5
$ws.Range("J20").Value = 1685.5385
$ws.Range("K20").Value = 1378.5
$ws.Range("L20").Value = 1685.5385
$ws.Range("M20").Value = -1131.5
$ws.Range("N20").Value = -2179.5385

$ws.Range("H134").Value = 2515.9343
$ws.Range("I134").Value = 2230.558
$ws.Range("J134").Value = 3197.6667
$ws.Range("K134").Value = 6691.674
$ws.Range("L134").Value = 9593.000100000001
$ws.Range("M134").Value = -4156.674
$ws.Range("N134").Value = -14663.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2603.868
$ws.Range("I31").Value = 1707.093
$ws.Range("J31").Value = 6460
$ws.Range("K31").Value = 1707.093
$ws.Range("L31").Value = 6460
$ws.Range("M31").Value = -1412.093
$ws.Range("N31").Value = -7050

$ws.Range("H34").Value = 2603.868
$ws.Range("I34").Value = 1707.093
$ws.Range("J34").Value = 6460
$ws.Range("K34").Value = 1707.093
$ws.Range("L34").Value = 6460
$ws.Range("M34").Value = -1505.093
$ws.Range("N34").Value = -6864

$ws.Range("H52").Value = 12800
$ws.Range("J52").Value = 12800
$ws.Range("L52").Value = 12800
$ws.Range("N52").Value = -13388

$ws.Range("H134").Value = 1480.4822
$ws.Range("I134").Value = 983.94446
$ws.Range("J134").Value = 2374.25
$ws.Range("K134").Value = 2951.83338
$ws.Range("L134").Value = 7122.75
$ws.Range("M134").Value = -416.83338
$ws.Range("N134").Value = -12192.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 478.6087
$ws.Range("I113").Value = 509.57144
$ws.Range("J113").Value = 465.0625
$ws.Range("K113").Value = 1528.71432
$ws.Range("L113").Value = 1395.1875
$ws.Range("M113").Value = 641.28568
$ws.Range("N113").Value = -5735.1875

$ws.Range("H122").Value = 869.56525
$ws.Range("J122").Value = 1142.8572
$ws.Range("L122").Value = 10285.7148
$ws.Range("N122").Value = -15185.7148

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4626.353
$ws.Range("I80").Value = 5510
$ws.Range("J80").Value = 3198.923
$ws.Range("K80").Value = 5510
$ws.Range("L80").Value = 3198.923
$ws.Range("M80").Value = -4512
$ws.Range("N80").Value = -5194.923

$ws.Range("H83").Value = 4626.353
$ws.Range("I83").Value = 5510
$ws.Range("J83").Value = 3198.923
$ws.Range("K83").Value = 27550
$ws.Range("L83").Value = 15994.615
$ws.Range("M83").Value = -22558
$ws.Range("N83").Value = -25978.615

$ws.Range("H122").Value = 2015.2858
$ws.Range("I122").Value = 1526.75
$ws.Range("J122").Value = 2666.6667
$ws.Range("K122").Value = 4580.25
$ws.Range("L122").Value = 8000.000100000001
$ws.Range("M122").Value = -2130.25
$ws.Range("N122").Value = -12900.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2617.5789
$ws.Range("I7").Value = 2599.8333
$ws.Range("K7").Value = 2599.8333
$ws.Range("M7").Value = -2487.8333

$ws.Range("H22").Value = 867.7692
$ws.Range("I22").Value = 870.1111
$ws.Range("J22").Value = 862.5
$ws.Range("K22").Value = 870.1111
$ws.Range("L22").Value = 862.5
$ws.Range("M22").Value = -575.1111
$ws.Range("N22").Value = -1452.5

$ws.Range("H27").Value = 867.7692
$ws.Range("I27").Value = 870.1111
$ws.Range("J27").Value = 862.5
$ws.Range("K27").Value = 870.1111
$ws.Range("L27").Value = 862.5
$ws.Range("M27").Value = -763.1111
$ws.Range("N27").Value = -1076.5

$ws.Range("H44").Value = 10000
$ws.Range("J44").Value = 10000
$ws.Range("L44").Value = 10000
$ws.Range("N44").Value = -10912

$ws.Range("H122").Value = 2340.1333
$ws.Range("I122").Value = 2333.3333
$ws.Range("J122").Value = 2350.3333
$ws.Range("K122").Value = 6999.999899999999
$ws.Range("L122").Value = 7050.999899999999
$ws.Range("M122").Value = -4549.999899999999
$ws.Range("N122").Value = -11950.9999

$ws.Range("H126").Value = 2617.5789
$ws.Range("I126").Value = 2599.8333
$ws.Range("K126").Value = 7799.499899999999
$ws.Range("M126").Value = -5329.499899999999

$ws.Range("H136").Value = 3550.8982
$ws.Range("I136").Value = 1978.05
$ws.Range("J136").Value = 6862.1577
$ws.Range("K136").Value = 5934.15
$ws.Range("L136").Value = 20586.4731
$ws.Range("M136").Value = -3384.15
$ws.Range("N136").Value = -25686.4731

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1720.6
$ws.Range("I122").Value = 1525.75
$ws.Range("K122").Value = 4577.25
$ws.Range("M122").Value = -2127.25

$ws.Range("H132").Value = 1644
$ws.Range("I132").Value = 1173.7097
$ws.Range("J132").Value = 2411.3157
$ws.Range("K132").Value = 3521.1291
$ws.Range("L132").Value = 7233.9471
$ws.Range("M132").Value = -991.1291000000001
$ws.Range("N132").Value = -12293.9471

$ws.Range("H136").Value = 832.09
$ws.Range("I136").Value = 537.0789
$ws.Range("J136").Value = 1766.2916
$ws.Range("K136").Value = 1611.2367
$ws.Range("L136").Value = 5298.8748
$ws.Range("M136").Value = 938.7633000000001
$ws.Range("N136").Value = -10398.8748
